$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.45
$ws.Range("H2").Value = 2.8
$ws.Range("I2").Value = 3.3
$ws.Range("L2").Value = 4.33
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("O2").Value = 1.67
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 3.1
$ws.Range("R2").Value = 1.36
$ws.Range("S2").Value = 1.67
$ws.Range("T2").Value = 2.1
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.53
$ws.Range("W2").Value = 5.5
$ws.Range("AC2").Value = 5
$ws.Range("AD2").Value = 5.5
$ws.Range("AE2").Value = 21
$ws.Range("AI2").Value = 15
$ws.Range("AK2").Value = 41
$ws.Range("AT2").Value = 2.1
$ws.Range("AU2").Value = 10
$ws.Range("AY2").Value = 41
$ws.Range("AZ2").Value = 81

# Row 3 updates
$ws.Range("N3").Value = 17
